# Updates cryptocurrency price (D) and 1h volume change (E) columns
# to match the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.612.33'
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").Value = '3.018.93'
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''380.01'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = '''102.48'
$ws.Range("E7").Value = '  +0.97%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''0.590'
$ws.Range("D10").Value = '''36.71'
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '''0.0860'
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").Value = '3.500.45'
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Value = '''18.42'
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = '''7.70'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").Value = '3.027.16'
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("D17").Value = '''0.974'
$ws.Range("E17").Value = '  -3.39%  '
$ws.Range("D18").Value = '''10.50'
$ws.Range("E18").Value = '  -14.48%  '
$ws.Range("D19").Value = '51.620.04'
$ws.Range("E19").Value = '  +1.28%  '
$ws.Range("D20").Value = '''3.07'
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").Value = '''12.40'
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").Value = '''69.94'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").Value = '''266.43'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  -7.19%  '
$ws.Range("D26").Value = '''8.20'
$ws.Range("E26").Value = '  +2.42%  '
$ws.Range("D27").Value = '''7.67'
$ws.Range("E27").Value = '  +11.31%  '
$ws.Range("D28").Value = '''0.171'
$ws.Range("E28").Value = '  +4.41%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '''26.16'
$ws.Range("E30").Value = '  +1.78%  '
$ws.Range("D31").Value = '''0.107'
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").Value = '''2.06'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").Value = '''50.40'
$ws.Range("D35").Value = '''33.70'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").Value = '''0.0446'
$ws.Range("E36").Value = '  +3.42%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '''3.29'
$ws.Range("E38").Value = '  +4.37%  '
$ws.Range("D39").Value = '''0.299'
$ws.Range("E39").Value = '  +16.99%  '
$ws.Range("D40").Value = '''16.96'
$ws.Range("E40").Value = '  +2.16%  '
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("D42").Value = '''126.99'
$ws.Range("E42").Value = '  +6.19%  '
$ws.Range("D44").Value = '''2.53'
$ws.Range("D45").Value = '''3.74'
$ws.Range("E45").Value = '  +6.10%  '
$ws.Range("D46").Value = '''21.51'
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("E47").Value = '  +3.91%  '
$ws.Range("D48").Value = '''2.39'
$ws.Range("E48").Value = '  +2.96%  '
$ws.Range("D49").Value = '2.025.53'
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").Value = '3.323.72'
$ws.Range("E50").Value = '  +2.39%  '
$ws.Range("D51").Value = '''0.0320'
$ws.Range("E51").Value = '  +0.89%  '
